$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    # Preserve the existing style while forcing the new value to be
    # stored as literal text (these cells hold numeric-looking strings
    # like "214.54" or "0.0500" that must NOT be coerced to numbers).
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "25.925.10"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.637.32"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.25%  "
Set-TextValue "D5" "214.54"
$ws.Range("E5").Value = "  -0.23%  "
Set-TextValue "D6" "0.507"
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  +0.53%  "
Set-TextValue "D10" "19.59"
$ws.Range("E10").Value = "  -0.58%  "
Set-TextValue "D11" "0.0795"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "1.863.81"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.24"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.632.12"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  -1.69%  "
$ws.Range("E16").Value = "  -0.67%  "
Set-TextValue "D17" "62.55"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "25.941.02"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E19").Value = "  +0.19%  "
Set-TextValue "D20" "193.82"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("E21").Value = "  -1.32%  "
Set-TextValue "D22" "9.91"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("E23").Value = "  -1.14%  "
Set-TextValue "D24" "1.82"
$ws.Range("E24").Value = "  +0.77%  "
Set-TextValue "D25" "143.73"
$ws.Range("E25").Value = "  +1.00%  "
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  +0.05%  "
Set-TextValue "D31" "0.0500"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("E34").Value = "  -2.47%  "
$ws.Range("E35").Value = "  +1.36%  "
Set-TextValue "D36" "0.902"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").Value = "1.137.84"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("E38").Value = "  -0.02%  "
Set-TextValue "D39" "2.46"
$ws.Range("E39").Value = "  -2.02%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("E41").Value = "  +0.10%  "
Set-TextValue "D42" "99.40"
$ws.Range("E42").Value = "  -1.42%  "
Set-TextValue "D43" "0.798"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("D45").Value = "1.772.79"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "0.0₆0113"
$ws.Range("E46").Value = "  +1.01%  "
Set-TextValue "D47" "56.44"
$ws.Range("E47").Value = "  +1.43%  "
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("E49").Value = "  -0.84%  "
Set-TextValue "D50" "7.63"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  -0.44%  "
